# Remove the now-unused "url.cell_type" column (column C) from the "survey" sheet.
# This shifts every column to its right one position to the left and, since the
# shared strings "url.cell_type" / "formula" become unreferenced, they drop out
# of the shared strings table automatically.
$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$survey.Columns.Item(3).Delete()

# The active/selected sheet moves from "initial" to "survey", with the
# selection on "survey" reset to B6 (and no special top-left scroll position).
$survey.Activate()
$survey.Range("B6").Select()
